$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.451.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.183.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.89"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.77%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -5.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0909"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.511.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.183.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("E17").Value = "  -5.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.363.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.42%  "
$ws.Range("E24").Value = "  -2.99%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.32%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0817"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.83%  "
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("E36").Value = "  -4.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.84%  "
$ws.Range("E42").Value = "  -7.97%  "
$ws.Range("E43").Value = "  -3.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "58.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.456"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("E50").Value = "  -2.57%  "
